# Applies the per-player stat corrections described in the commit
# ("Update analysis scripts, fix Ligue 1 aliases, and update prediction
# data") to the Hamburger SV player-stats sheet. Every cell below is a
# plain numeric value (no formulas), so we just overwrite Range.Value
# for each touched cell, grouped by the player row it belongs to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AA2").Value = 0.30353649
$ws.Range("AF2").Value = 61.538461538462
$ws.Range("AG2").Value = 78
$ws.Range("BT2").Value = 30
$ws.Range("DG2").Value = 58

# Row 3
$ws.Range("E3").Value = 6.6333333333333
$ws.Range("K3").Value = 6.0372
$ws.Range("BB3").Value = 39.230769230769
$ws.Range("BF3").Value = 44.444444444444
$ws.Range("BG3").Value = 142
$ws.Range("BQ3").Value = 139.3
$ws.Range("CC3").Value = 9
$ws.Range("CD3").Value = 6
$ws.Range("CU3").Value = 79
$ws.Range("CV3").Value = 35

# Row 4
$ws.Range("AF4").Value = 70
$ws.Range("AG4").Value = 50
$ws.Range("BB4").Value = 30.555555555556
$ws.Range("BF4").Value = 28.571428571429
$ws.Range("BT4").Value = 15
$ws.Range("CU4").Value = 25
$ws.Range("CV4").Value = 15
$ws.Range("DG4").Value = 34

# Row 10
$ws.Range("E10").Value = 7.0375
$ws.Range("K10").Value = 2.9036
$ws.Range("AA10").Value = 2.30480443
$ws.Range("AB10").Value = 657
$ws.Range("AK10").Value = 35
$ws.Range("AL10").Value = 58.333333333333
$ws.Range("AY10").Value = 16
$ws.Range("AZ10").Value = 53.333333333333
$ws.Range("BA10").Value = 46
$ws.Range("BB10").Value = 50
$ws.Range("BC10").Value = 44
$ws.Range("BD10").Value = 56.410256410256
$ws.Range("BF10").Value = 14.285714285714
$ws.Range("BQ10").Value = 112.6
$ws.Range("CC10").Value = 45
$ws.Range("CS10").Value = 30
$ws.Range("CU10").Value = 46
$ws.Range("CV10").Value = 12
$ws.Range("CZ10").Value = 60

# Row 12
$ws.Range("E12").Value = 6.7666666666667
$ws.Range("AB12").Value = 859
$ws.Range("AO12").Value = 32
$ws.Range("AY12").Value = 5
$ws.Range("AZ12").Value = 38.461538461538
$ws.Range("BA12").Value = 112
$ws.Range("BB12").Value = 58.947368421053
$ws.Range("BC12").Value = 89
$ws.Range("BD12").Value = 62.237762237762
$ws.Range("BQ12").Value = 121.8
$ws.Range("CS12").Value = 13

# Row 15
$ws.Range("AZ15").Value = 50
$ws.Range("CB15").Value = 5
$ws.Range("CC15").Value = 5
$ws.Range("CD15").Value = 3
$ws.Range("CS15").Value = 6

# Row 16
$ws.Range("AA16").Value = 0.97928686
$ws.Range("AE16").Value = 801
$ws.Range("AF16").Value = 90.406320541761
$ws.Range("AG16").Value = 886
$ws.Range("AH16").Value = 512
$ws.Range("AI16").Value = 289
$ws.Range("AO16").Value = 27
$ws.Range("AT16").Value = 69
$ws.Range("BG16").Value = 116
$ws.Range("BT16").Value = 85
$ws.Range("CC16").Value = 46
$ws.Range("CD16").Value = 35
$ws.Range("DA16").Value = 17
$ws.Range("DB16").Value = 41.463414634146
$ws.Range("DF16").Value = 545
$ws.Range("DG16").Value = 341

# Row 17
$ws.Range("AA17").Value = 0.58432916
$ws.Range("AH17").Value = 694
$ws.Range("AI17").Value = 281
$ws.Range("BB17").Value = 68.075117370892
$ws.Range("BF17").Value = 76.642335766423
$ws.Range("CU17").Value = 68
$ws.Range("CV17").Value = 32
$ws.Range("DF17").Value = 735
$ws.Range("DG17").Value = 424

# Row 18
$ws.Range("AA18").Value = 0.42965376
$ws.Range("AE18").Value = 154
$ws.Range("AF18").Value = 78.571428571429
$ws.Range("AG18").Value = 196
$ws.Range("AH18").Value = 87
$ws.Range("DF18").Value = 101

# Row 19
$ws.Range("K19").Value = 1.0858
$ws.Range("AT19").Value = 44
$ws.Range("BG19").Value = 359

# Row 20
$ws.Range("AA20").Value = 0.60811991
$ws.Range("AB20").Value = 832
$ws.Range("AE20").Value = 541
$ws.Range("AF20").Value = 84.929356357928
$ws.Range("AG20").Value = 637
$ws.Range("AH20").Value = 394
$ws.Range("AK20").Value = 24
$ws.Range("AL20").Value = 36.363636363636
$ws.Range("AO20").Value = 15
$ws.Range("AT20").Value = 83
$ws.Range("BB20").Value = 66.666666666667
$ws.Range("BF20").Value = 68.181818181818
$ws.Range("CC20").Value = 54
$ws.Range("CD20").Value = 29
$ws.Range("CU20").Value = 33
$ws.Range("CV20").Value = 14
$ws.Range("CZ20").Value = 66
$ws.Range("DF20").Value = 427

# Row 23
$ws.Range("E23").Value = 6.4
$ws.Range("AW23").Value = 0
$ws.Range("BQ23").Value = 89.59999999999999
